$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.489.81"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").Value = "1.724.28"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5379"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07726"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.637"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.726.73"
$ws.Range("E13").Value = "  +3.70%  "
$ws.Range("D14").Value = "1.960.83"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5878"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("D16").Value = "0.0₅8279"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "27.506.19"
$ws.Range("E18").Value = "  +5.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.84%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.736"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.101"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.689"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.48%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1231"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.412"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05557"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.549"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.469"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.658"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9594"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.812"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5922"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01643"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.863"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8561"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").Value = "1.055.44"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "1.866.95"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +10.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4439"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05275"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.32%  "
